$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5389986038208008
$ws.Range("E2").Value = 126.6426861646778
$ws.Range("F2").Value = 0.004364484159258046
$ws.Range("G2").Value = 0.003685497530278848
$ws.Range("H2").Value = 0.003281981489578283
$ws.Range("I2").Value = 0.003069505575590013
$ws.Range("J2").Value = 0.003069505575590013
$ws.Range("K2").Value = 0.003069505575590013
$ws.Range("L2").Value = 0.003006090023300911
$ws.Range("M2").Value = 0.002824927465330597
$ws.Range("N2").Value = 0.002635798222585095
$ws.Range("O2").Value = 0.002635798222585095
$ws.Range("P2").Value = 0.002635798222585095
$ws.Range("Q2").Value = 0.002635798222585095
$ws.Range("R2").Value = 0.002635798222585095
$ws.Range("S2").Value = 0.002635798222585095
$ws.Range("T2").Value = 0.002604703742776248
$ws.Range("U2").Value = 0.002604703742776248
$ws.Range("V2").Value = 0.002583682949108359
$ws.Range("W2").Value = 0.002540257066612453
$ws.Range("X2").Value = 0.002505728536850786
$ws.Range("Y2").Value = 0.002468668346290015

$ws.Range("C3").Value = 0.449979305267334
$ws.Range("E3").Value = 114.0170331663157
$ws.Range("F3").Value = 0.004448212927907743
$ws.Range("G3").Value = 0.003414919473505197
$ws.Range("H3").Value = 0.003188512887160494
$ws.Range("I3").Value = 0.003188512887160494
$ws.Range("J3").Value = 0.003028704169767002
$ws.Range("K3").Value = 0.002834342829621419
$ws.Range("L3").Value = 0.002730174008564445
$ws.Range("M3").Value = 0.002686869856043151
$ws.Range("N3").Value = 0.002452970819439228
$ws.Range("O3").Value = 0.002452970819439228
$ws.Range("P3").Value = 0.002452970819439228
$ws.Range("Q3").Value = 0.002452970819439228
$ws.Range("R3").Value = 0.002452970819439228
$ws.Range("S3").Value = 0.002452970819439228
$ws.Range("T3").Value = 0.00237631774161823
$ws.Range("U3").Value = 0.002348547453875476
$ws.Range("V3").Value = 0.002315017083751369
$ws.Range("W3").Value = 0.002250480270419383
$ws.Range("X3").Value = 0.002250480270419383
$ws.Range("Y3").Value = 0.002222554252754692

$ws.Range("C4").Value = 0.4610099792480469
$ws.Range("E4").Value = 114.1786462306791
$ws.Range("F4").Value = 0.00431930524465792
$ws.Range("G4").Value = 0.003686233814105452
$ws.Range("H4").Value = 0.003417906352002872
$ws.Range("I4").Value = 0.003193814751914287
$ws.Range("J4").Value = 0.002981673028344625
$ws.Range("K4").Value = 0.002926863113464181
$ws.Range("L4").Value = 0.00262925099278854
$ws.Range("M4").Value = 0.00262925099278854
$ws.Range("N4").Value = 0.00259476189406099
$ws.Range("O4").Value = 0.002591766419197783
$ws.Range("P4").Value = 0.002591766419197783
$ws.Range("Q4").Value = 0.002491120409592935
$ws.Range("R4").Value = 0.002301550129094868
$ws.Range("S4").Value = 0.002301550129094868
$ws.Range("T4").Value = 0.002301550129094868
$ws.Range("U4").Value = 0.002301550129094868
$ws.Range("V4").Value = 0.002265601786192823
$ws.Range("W4").Value = 0.002265601786192823
$ws.Range("X4").Value = 0.002256498845168114
$ws.Range("Y4").Value = 0.002225704604886531

$ws.Range("C5").Value = 0.4689986705780029
$ws.Range("E5").Value = 111.8916679397062
$ws.Range("F5").Value = 0.004350581652676185
$ws.Range("G5").Value = 0.003576971710480415
$ws.Range("H5").Value = 0.003513894384160952
$ws.Range("I5").Value = 0.003210117945853883
$ws.Range("J5").Value = 0.003078236728401711
$ws.Range("K5").Value = 0.002962184825010395
$ws.Range("L5").Value = 0.002757360494654961
$ws.Range("M5").Value = 0.002627148751013137
$ws.Range("N5").Value = 0.002451135001929021
$ws.Range("O5").Value = 0.002451135001929021
$ws.Range("P5").Value = 0.002451135001929021
$ws.Range("Q5").Value = 0.002451135001929021
$ws.Range("R5").Value = 0.002421197347380283
$ws.Range("S5").Value = 0.002366126374177894
$ws.Range("T5").Value = 0.002336921034986967
$ws.Range("U5").Value = 0.0022379755790026
$ws.Range("V5").Value = 0.002224648279243588
$ws.Range("W5").Value = 0.002181124131378287
$ws.Range("X5").Value = 0.002181124131378287
$ws.Range("Y5").Value = 0.002181124131378287

$ws.Range("C6").Value = 0.5079703330993652
$ws.Range("E6").Value = 116.3791679893038
$ws.Range("F6").Value = 0.004254990022205075
$ws.Range("G6").Value = 0.003667130267142701
$ws.Range("H6").Value = 0.00317091710152404
$ws.Range("I6").Value = 0.002996337006444205
$ws.Range("J6").Value = 0.002969378274388894
$ws.Range("K6").Value = 0.002854746259816674
$ws.Range("L6").Value = 0.002835007497225079
$ws.Range("M6").Value = 0.002835007497225079
$ws.Range("N6").Value = 0.002719192874422387
$ws.Range("O6").Value = 0.002715002246197371
$ws.Range("P6").Value = 0.002641810965752478
$ws.Range("Q6").Value = 0.002641810965752478
$ws.Range("R6").Value = 0.00263668435348508
$ws.Range("S6").Value = 0.002541752682698157
$ws.Range("T6").Value = 0.002395247034280091
$ws.Range("U6").Value = 0.002362522109695376
$ws.Range("V6").Value = 0.002314518946876165
$ws.Range("W6").Value = 0.00231125733281932
$ws.Range("X6").Value = 0.002303139617990885
$ws.Range("Y6").Value = 0.002268599765873369

$ws.Range("C7").Value = 0.607003927230835
$ws.Range("E7").Value = 111.8859606275564
$ws.Range("F7").Value = 0.00438295572147182
$ws.Range("G7").Value = 0.003344109421014846
$ws.Range("H7").Value = 0.003273751451447468
$ws.Range("I7").Value = 0.003195090160578947
$ws.Range("J7").Value = 0.00307404770036543
$ws.Range("K7").Value = 0.002913492033346406
$ws.Range("L7").Value = 0.002783756172072764
$ws.Range("M7").Value = 0.002460185038568832
$ws.Range("N7").Value = 0.002460185038568832
$ws.Range("O7").Value = 0.002460185038568832
$ws.Range("P7").Value = 0.002460185038568832
$ws.Range("Q7").Value = 0.002460185038568832
$ws.Range("R7").Value = 0.002412689977614047
$ws.Range("S7").Value = 0.002368030305463471
$ws.Range("T7").Value = 0.002276445955006187
$ws.Range("U7").Value = 0.002232737977554432
$ws.Range("V7").Value = 0.002232737977554432
$ws.Range("W7").Value = 0.002217768294990465
$ws.Range("X7").Value = 0.00219548548218541
$ws.Range("Y7").Value = 0.002181012877730144

$ws.Range("C8").Value = 0.6359982490539551
$ws.Range("E8").Value = 114.0765903044612
$ws.Range("F8").Value = 0.00432329529302083
$ws.Range("G8").Value = 0.003552023658856618
$ws.Range("H8").Value = 0.003290769825804059
$ws.Range("I8").Value = 0.003218745790353153
$ws.Range("J8").Value = 0.003081407328576971
$ws.Range("K8").Value = 0.002799039809469725
$ws.Range("L8").Value = 0.002722081274767769
$ws.Range("M8").Value = 0.002627808956216714
$ws.Range("N8").Value = 0.002599128195008595
$ws.Range("O8").Value = 0.002503588136091723
$ws.Range("P8").Value = 0.002503588136091723
$ws.Range("Q8").Value = 0.002485952459278672
$ws.Range("R8").Value = 0.002448559044323967
$ws.Range("S8").Value = 0.002339981822858204
$ws.Range("T8").Value = 0.002339981822858204
$ws.Range("U8").Value = 0.002309790881299633
$ws.Range("V8").Value = 0.002297559184249305
$ws.Range("W8").Value = 0.002291856549304691
$ws.Range("X8").Value = 0.002237755450708623
$ws.Range("Y8").Value = 0.002223715210613278

$ws.Range("C9").Value = 0.6570007801055908
$ws.Range("E9").Value = 119.2989798631188
$ws.Range("F9").Value = 0.004448212927907743
$ws.Range("G9").Value = 0.003584895951310413
$ws.Range("H9").Value = 0.003210983300468407
$ws.Range("I9").Value = 0.003031272813248038
$ws.Range("J9").Value = 0.002960049507930509
$ws.Range("K9").Value = 0.00292400691432459
$ws.Range("L9").Value = 0.002754679555196702
$ws.Range("M9").Value = 0.00275339199792444
$ws.Range("N9").Value = 0.00253201840521467
$ws.Range("O9").Value = 0.00253201840521467
$ws.Range("P9").Value = 0.00253201840521467
$ws.Range("Q9").Value = 0.00253201840521467
$ws.Range("R9").Value = 0.00253201840521467
$ws.Range("S9").Value = 0.002509588556658429
$ws.Range("T9").Value = 0.002509588556658429
$ws.Range("U9").Value = 0.002459512536976322
$ws.Range("V9").Value = 0.002459512536976322
$ws.Range("W9").Value = 0.002393320314703705
$ws.Range("X9").Value = 0.002346985773831007
$ws.Range("Y9").Value = 0.002325516176668982

$ws.Range("C10").Value = 0.5919914245605469
$ws.Range("E10").Value = 115.0709959736505
$ws.Range("F10").Value = 0.004264709933331066
$ws.Range("G10").Value = 0.003431084509297379
$ws.Range("H10").Value = 0.0033446797172811
$ws.Range("I10").Value = 0.003141187245167948
$ws.Range("J10").Value = 0.002911726273634722
$ws.Range("K10").Value = 0.002582346345895457
$ws.Range("L10").Value = 0.002582346345895457
$ws.Range("M10").Value = 0.002582346345895457
$ws.Range("N10").Value = 0.002582346345895457
$ws.Range("O10").Value = 0.002582346345895457
$ws.Range("P10").Value = 0.002345909127039388
$ws.Range("Q10").Value = 0.002345909127039388
$ws.Range("R10").Value = 0.002345909127039388
$ws.Range("S10").Value = 0.002345909127039388
$ws.Range("T10").Value = 0.002345909127039388
$ws.Range("U10").Value = 0.002289833849279365
$ws.Range("V10").Value = 0.002289833849279365
$ws.Range("W10").Value = 0.002256956399206921
$ws.Range("X10").Value = 0.002256956399206921
$ws.Range("Y10").Value = 0.002243099336718333

$ws.Range("C11").Value = 0.4180037975311279
$ws.Range("E11").Value = 114.3440445586621
$ws.Range("F11").Value = 0.004279189778590753
$ws.Range("G11").Value = 0.003445253278320225
$ws.Range("H11").Value = 0.003385610598878073
$ws.Range("I11").Value = 0.00305647854442224
$ws.Range("J11").Value = 0.002694268706645441
$ws.Range("K11").Value = 0.002694268706645441
$ws.Range("L11").Value = 0.002689140058575563
$ws.Range("M11").Value = 0.002689140058575563
$ws.Range("N11").Value = 0.002679905176483641
$ws.Range("O11").Value = 0.002529912645824656
$ws.Range("P11").Value = 0.00242396822757141
$ws.Range("Q11").Value = 0.00242396822757141
$ws.Range("R11").Value = 0.00242396822757141
$ws.Range("S11").Value = 0.002415648182219303
$ws.Range("T11").Value = 0.002369266043796455
$ws.Range("U11").Value = 0.002333234887971618
$ws.Range("V11").Value = 0.00228368664805378
$ws.Range("W11").Value = 0.002262307570470131
$ws.Range("X11").Value = 0.002251117901553646
$ws.Range("Y11").Value = 0.002228928743833568
